$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.011499999999995
$ws.Range("E4").Value = 12.6196
$ws.Range("A9").Value = -20.23869999999999
$ws.Range("E10").Value = 12.4713
$ws.Range("A18").Value = -23.03790000000001
$ws.Range("A20").Value = -22.13660000000003
$ws.Range("C21").Value = -12.90300000000001
